$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E13")
foreach ($cell in $rng.Cells) {
    $cell.Value2 = [math]::Round([double]$cell.Value2, 0)
}
